$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.942.41"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "2.909.76"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.06"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.47"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "2.908.83"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.80"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.58"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "3.390.31"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("D17").Value = "60.871.37"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.73"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").Value = "2.907.49"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.70"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.52"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.684"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.09"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.00"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.92"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.10"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.33"
$ws.Range("E29").Value = "  +5.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.07"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.63"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "0.0₃0855"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.10"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.62"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.66"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.125"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.63"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.289"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.17"
$ws.Range("E44").Value = "  -5.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "378.21"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0349"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").Value = "2.709.03"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.77"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.09"
$ws.Range("E50").Value = "  -6.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.107"
$ws.Range("E51").Value = "  +0.49%  "
